$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Trae Young"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Atlanta Hawks"
$ws.Range("A3").Value = "Caris LeVert"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Atlanta Hawks"
$ws.Range("A4").Value = "Scoot Henderson"
$ws.Range("B4").Value = "PG"
$ws.Range("C4").Value = "Portland Trail Blazers"
$ws.Range("A5").Value = "Desmond Bane"
$ws.Range("B5").Value = "SG,SF"
$ws.Range("C5").Value = "Memphis Grizzlies"
$ws.Range("A6").Value = "Jalen Williams"
$ws.Range("B6").Value = "SG,SF,PF,C"
$ws.Range("C6").Value = "Oklahoma City Thunder"
$ws.Range("A7").Value = "Jabari Smith Jr."
$ws.Range("B7").Value = "PF,C"
$ws.Range("C7").Value = "Houston Rockets"
$ws.Range("A8").Value = "Walker Kessler"
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "Utah Jazz"
$ws.Range("A9").Value = "Jonas Valanciunas"
$ws.Range("B9").Value = "C"
$ws.Range("C9").Value = "Sacramento Kings"
$ws.Range("A10").Value = "Myles Turner"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Indiana Pacers"
$ws.Range("A11").Value = "Kawhi Leonard"
$ws.Range("B11").Value = "SG,SF,PF"
$ws.Range("C11").Value = "LA Clippers"
$ws.Range("A12").Value = "LeBron James"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "Los Angeles Lakers"
$ws.Range("A13").Value = "Jalen Brunson"
$ws.Range("B13").Value = "PG"
$ws.Range("C13").Value = "New York Knicks"
$ws.Range("A14").Value = "Devin Booker"
$ws.Range("B14").Value = "PG,SG"
$ws.Range("C14").Value = "Phoenix Suns"
$ws.Range("A15").Value = "Immanuel Quickley"
$ws.Range("B15").Value = "PG,SG"
$ws.Range("C15").Value = "Toronto Raptors"
$ws.Range("A16").Value = "Nicolas Claxton"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Brooklyn Nets"
$ws.Range("A17").Value = "Brandon Ingram"
$ws.Range("B17").Value = "SG,SF,PF"
$ws.Range("C17").Value = "Toronto Raptors"
$ws.Range("A18").Value = "Norman Powell"
$ws.Range("B18").Value = "SG,SF"
$ws.Range("C18").Value = "LA Clippers"
$ws.Range("A19").Value = "Jimmy Butler III"
$ws.Range("B19").Value = "SF,PF"
$ws.Range("C19").Value = "Golden State Warriors"
